$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style  # reference: an untouched, default-styled cell

$ws.Range("D2").Value = '60.302.38'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '3.301.60'
$ws.Range("E3").Value = '  -3.46%  '
$ws.Range("D5").Value = '''557.10'
$ws.Range("E5").Value = '  -3.85%  '
$ws.Range("D6").Value = '''141.31'
$ws.Range("E6").Value = '  -8.49%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.302.29'
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("E9").Value = '  -3.62%  '
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").Value = '  -5.00%  '
$ws.Range("E12").Value = '  -2.72%  '
$ws.Range("D13").Value = '3.865.45'
$ws.Range("E13").Value = '  -3.44%  '
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '''26.66'
$ws.Range("E15").Value = '  -7.01%  '
$ws.Range("D16").Value = '3.302.13'
$ws.Range("E16").Value = '  -3.41%  '
$ws.Range("E17").Value = '  -4.96%  '
$ws.Range("D18").Value = '60.289.31'
$ws.Range("E18").Value = '  -2.80%  '
$ws.Range("E19").Value = '  -7.83%  '
$ws.Range("D20").Value = '''13.73'
$ws.Range("E20").Value = '  -4.99%  '
$ws.Range("D21").Value = '''8.53'
$ws.Range("E21").Value = '  -4.94%  '
$ws.Range("D22").Value = '''373.84'
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '''72.28'
$ws.Range("E24").Value = '  -4.94%  '
$ws.Range("D25").Value = '''0.532'
$ws.Range("E25").Value = '  -6.47%  '
$ws.Range("D26").Value = '3.432.66'
$ws.Range("E26").Value = '  -3.61%  '
$ws.Range("D27").Value = '''0.0000102'
$ws.Range("E27").Value = '  -9.32%  '
$ws.Range("D28").Value = '''0.174'
$ws.Range("E28").Value = '  -2.30%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '''7.04'
$ws.Range("E30").Value = '  -8.12%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  -5.03%  '
$ws.Range("D33").Value = '''7.41'
$ws.Range("E33").Value = '  -6.06%  '
$ws.Range("D34").Value = '''22.57'
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("D35").Value = '''1.25'
$ws.Range("E35").Value = '  -6.02%  '
$ws.Range("E36").Value = '  -9.76%  '
$ws.Range("D37").Value = '''165.73'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("E38").Value = '  -4.63%  '
$ws.Range("E39").Value = '  -4.98%  '
$ws.Range("D41").Value = '''0.0723'
$ws.Range("E41").Value = '  -7.85%  '
$ws.Range("D42").Value = '''25.48'
$ws.Range("E42").Value = '  -17.87%  '
$ws.Range("D43").Value = '''41.66'
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("D45").Value = '''1.13'
$ws.Range("E45").Value = '  -3.71%  '
$ws.Range("E46").Value = '  -7.62%  '
$ws.Range("E47").Value = '  -6.34%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").Value = '2.323.53'
$ws.Range("E49").Value = '  -8.95%  '
$ws.Range("D50").Value = '''21.59'
$ws.Range("E50").Value = '  -7.17%  '
$ws.Range("D51").Value = '''6.34'
$ws.Range("E51").Value = '  -6.94%  '

# The apostrophe prefix marks those cells with a "quote prefix" style; restore the
# original (default) cell style so only the cell VALUE changed, matching the source file.
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").Style = $plainStyle
$ws.Range("D15").Style = $plainStyle
$ws.Range("D20").Style = $plainStyle
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").Style = $plainStyle
$ws.Range("D24").Style = $plainStyle
$ws.Range("D25").Style = $plainStyle
$ws.Range("D27").Style = $plainStyle
$ws.Range("D28").Style = $plainStyle
$ws.Range("D29").Style = $plainStyle
$ws.Range("D30").Style = $plainStyle
$ws.Range("D33").Style = $plainStyle
$ws.Range("D34").Style = $plainStyle
$ws.Range("D35").Style = $plainStyle
$ws.Range("D37").Style = $plainStyle
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").Style = $plainStyle
$ws.Range("D43").Style = $plainStyle
$ws.Range("D45").Style = $plainStyle
$ws.Range("D50").Style = $plainStyle
$ws.Range("D51").Style = $plainStyle
